# Scheduled runner refresh: pulls latest Universalis market-board averages and
# re-derives the LevePriceNQ/HQ + LeveProfitNQ/HQ columns for every affected leve row.
# (currentAveragePrice[, NQ, HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ] -- columns H..N)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 41874.2
$ws.Range("J3").Value = 41874.2
$ws.Range("L3").Value = 41874.2
$ws.Range("N3").Value = -42102.2
# Row 17
$ws.Range("H17").Value = 827.6
$ws.Range("J17").Value = 827.6
$ws.Range("L17").Value = 2482.8
$ws.Range("N17").Value = -2818.8
# Row 21
$ws.Range("H21").Value = 72514.75
$ws.Range("I21").Value = 70019
$ws.Range("J21").Value = 75010.5
$ws.Range("K21").Value = 70019
$ws.Range("L21").Value = 75010.5
$ws.Range("M21").Value = -69551
$ws.Range("N21").Value = -75946.5
# Row 23
$ws.Range("H23").Value = 72514.75
$ws.Range("I23").Value = 70019
$ws.Range("J23").Value = 75010.5
$ws.Range("K23").Value = 70019
$ws.Range("L23").Value = 75010.5
$ws.Range("M23").Value = -69785
$ws.Range("N23").Value = -75478.5
# Row 34
$ws.Range("H34").Value = 3065
$ws.Range("I34").Value = 2916.1538
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 2916.1538
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -2713.1538
$ws.Range("N34").Value = -5406
# Row 36
$ws.Range("H36").Value = 3065
$ws.Range("I36").Value = 2916.1538
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 2916.1538
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -2201.1538
$ws.Range("N36").Value = -6430
# Row 102
$ws.Range("H102").Value = 41874.2
$ws.Range("J102").Value = 41874.2
$ws.Range("L102").Value = 41874.2
$ws.Range("N102").Value = -48364.2
# Row 103
$ws.Range("H103").Value = 2112.5
$ws.Range("I103").Value = 900
$ws.Range("K103").Value = 2700
$ws.Range("M103").Value = -2114
# Row 137
$ws.Range("H137").Value = 2176.5
$ws.Range("I137").Value = 2277.4
$ws.Range("K137").Value = 6832.200000000001
$ws.Range("M137").Value = -4282.200000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5287.725
$ws.Range("I32").Value = 4743.108
$ws.Range("K32").Value = 4743.108
$ws.Range("M32").Value = -4456.108
# Row 92
$ws.Range("H92").Value = 31840
$ws.Range("J92").Value = 31840
$ws.Range("L92").Value = 31840
$ws.Range("N92").Value = -36832
# Row 95
$ws.Range("H95").Value = 19047
$ws.Range("J95").Value = 19047
$ws.Range("L95").Value = 19047
$ws.Range("N95").Value = -24539

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 95
$ws.Range("H95").Value = 12574
$ws.Range("J95").Value = 12574
$ws.Range("L95").Value = 12574
$ws.Range("N95").Value = -18066

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 10755377
$ws.Range("I31").Value = 1891.9445
$ws.Range("J31").Value = 25644818
$ws.Range("K31").Value = 1891.9445
$ws.Range("L31").Value = 25644818
$ws.Range("M31").Value = -1596.9445
$ws.Range("N31").Value = -25645408
# Row 34
$ws.Range("H34").Value = 10755377
$ws.Range("I34").Value = 1891.9445
$ws.Range("J34").Value = 25644818
$ws.Range("K34").Value = 1891.9445
$ws.Range("L34").Value = 25644818
$ws.Range("M34").Value = -1689.9445
$ws.Range("N34").Value = -25645222
# Row 134
$ws.Range("H134").Value = 1591.238
$ws.Range("I134").Value = 990.7368
$ws.Range("J134").Value = 2087.3044
$ws.Range("K134").Value = 2972.2104
$ws.Range("L134").Value = 6261.9132
$ws.Range("M134").Value = -437.2103999999999
$ws.Range("N134").Value = -11331.9132

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 375.4375
$ws.Range("I2").Value = 45
$ws.Range("J2").Value = 705.875
$ws.Range("K2").Value = 270
$ws.Range("L2").Value = 4235.25
$ws.Range("M2").Value = -157
$ws.Range("N2").Value = -4461.25
# Row 9
$ws.Range("H9").Value = 2980
$ws.Range("J9").Value = 2980
$ws.Range("L9").Value = 8940
$ws.Range("N9").Value = -9388
# Row 15
$ws.Range("H15").Value = 497.8
$ws.Range("I15").Value = 295.6
$ws.Range("J15").Value = 700
$ws.Range("K15").Value = 886.8000000000001
$ws.Range("L15").Value = 2100
$ws.Range("M15").Value = -746.8000000000001
$ws.Range("N15").Value = -2380
# Row 22
$ws.Range("H22").Value = 1741.5714
$ws.Range("I22").Value = 2792.75
$ws.Range("J22").Value = 340
$ws.Range("K22").Value = 8378.25
$ws.Range("L22").Value = 1020
$ws.Range("M22").Value = -8209.25
$ws.Range("N22").Value = -1358
# Row 27
$ws.Range("H27").Value = 1741.5714
$ws.Range("I27").Value = 2792.75
$ws.Range("J27").Value = 340
$ws.Range("K27").Value = 8378.25
$ws.Range("L27").Value = 1020
$ws.Range("M27").Value = -8276.25
$ws.Range("N27").Value = -1224
# Row 69
$ws.Range("H69").Value = 2818.5454
$ws.Range("I69").Value = 1001.3333
$ws.Range("K69").Value = 3003.9999
$ws.Range("M69").Value = -2192.9999
# Row 72
$ws.Range("H72").Value = 2818.5454
$ws.Range("I72").Value = 1001.3333
$ws.Range("K72").Value = 9011.9997
$ws.Range("M72").Value = -4955.9997
# Row 92
$ws.Range("H92").Value = 1164.2
$ws.Range("J92").Value = 1382.8182
$ws.Range("L92").Value = 4148.4546
$ws.Range("N92").Value = -6644.4546
# Row 99
$ws.Range("H99").Value = 2030.3
$ws.Range("I99").Value = 1168.75
$ws.Range("J99").Value = 2604.6667
$ws.Range("K99").Value = 3506.25
$ws.Range("L99").Value = 7814.000100000001
$ws.Range("M99").Value = -1260.25
$ws.Range("N99").Value = -12306.0001
# Row 107
$ws.Range("H107").Value = 884
$ws.Range("I107").Value = 362
$ws.Range("J107").Value = 1145
$ws.Range("K107").Value = 1086
$ws.Range("L107").Value = 3435
$ws.Range("M107").Value = 834
$ws.Range("N107").Value = -7275
# Row 133
$ws.Range("H133").Value = 3841.25
$ws.Range("J133").Value = 4960
$ws.Range("L133").Value = 14880
$ws.Range("N133").Value = -25000
# Row 137
$ws.Range("H137").Value = 5514.8184
$ws.Range("I137").Value = 2926
$ws.Range("J137").Value = 7672.1665
$ws.Range("K137").Value = 8778
$ws.Range("L137").Value = 23016.4995
$ws.Range("M137").Value = -3678
$ws.Range("N137").Value = -33216.49950000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 8000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15302
# Row 46
$ws.Range("H46").Value = 29998
$ws.Range("J46").Value = 29998
$ws.Range("L46").Value = 29998
$ws.Range("N46").Value = -30310
# Row 57
$ws.Range("H57").Value = 17749.75
$ws.Range("J57").Value = 20333
$ws.Range("L57").Value = 20333
$ws.Range("N57").Value = -21973
# Row 80
$ws.Range("H80").Value = 2818.6365
$ws.Range("I80").Value = 2517.5
$ws.Range("J80").Value = 3180
$ws.Range("K80").Value = 2517.5
$ws.Range("L80").Value = 3180
$ws.Range("M80").Value = -1519.5
$ws.Range("N80").Value = -5176
# Row 83
$ws.Range("H83").Value = 2818.6365
$ws.Range("I83").Value = 2517.5
$ws.Range("J83").Value = 3180
$ws.Range("K83").Value = 12587.5
$ws.Range("L83").Value = 15900
$ws.Range("M83").Value = -7595.5
$ws.Range("N83").Value = -25884
# Row 95
$ws.Range("H95").Value = 10581
$ws.Range("J95").Value = 10581
$ws.Range("L95").Value = 10581
$ws.Range("N95").Value = -16073

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 94
$ws.Range("H94").Value = 29500
$ws.Range("J94").Value = 29500
$ws.Range("L94").Value = 29500
$ws.Range("N94").Value = -30852
# Row 105
$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988
# Row 132
$ws.Range("H132").Value = 58827150
$ws.Range("I132").Value = 71431896
$ws.Range("K132").Value = 214295688
$ws.Range("M132").Value = -214293158

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
# Row 98
$ws.Range("H98").Value = 38863.332
$ws.Range("J98").Value = 38863.332
$ws.Range("L98").Value = 38863.332
$ws.Range("N98").Value = -44853.332
# Row 105
$ws.Range("H105").Value = 23538.334
$ws.Range("J105").Value = 23538.334
$ws.Range("L105").Value = 23538.334
$ws.Range("N105").Value = -30526.334
# Row 141
$ws.Range("H141").Value = 35666.668
$ws.Range("J141").Value = 35666.668
$ws.Range("L141").Value = 35666.668
$ws.Range("N141").Value = -46026.668

